# Refresh cryptos price list: updates Price (D) and Volume(1h) (E) columns with
# newly fetched figures. Also MXToken (row 38) overtakes HuobiToken (row 37) in
# ranking, so rows 37/38 swap their Coin/Link/Price/Volume contents.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.650.75'
$ws.Range("E2").Value = '  +0.77%  '
$ws.Range("D3").Value = '1.565.01'
$ws.Range("E3").Value = '  -0.57%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.80%  '
$ws.Range("E6").Value = '  -0.64%  '
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '25.02'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.49%  '
$ws.Range("E9").Value = '  -0.75%  '
$ws.Range("E10").Value = '  -0.30%  '
$ws.Range("E11").Value = '  +0.02%  '
$ws.Range("D12").Value = '1.788.87'
$ws.Range("E12").Value = '  -0.57%  '
$ws.Range("D13").Value = '1.560.26'
$ws.Range("E13").Value = '  -0.82%  '
$ws.Range("D14").Value = '28.658.27'
$ws.Range("E14").Value = '  +0.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.513'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.45%  '
$ws.Range("E17").Value = '  -0.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '231.42'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.37'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.62%  '
$ws.Range("E22").Value = '  -1.24%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.96'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.20'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.77'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.104'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.11%  '
$ws.Range("E28").Value = '  -0.15%  '
$ws.Range("E30").Value = '  -4.71%  '
$ws.Range("E31").Value = '  -1.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.17'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.78%  '
$ws.Range("D33").Value = '1.391.09'
$ws.Range("E33").Value = '  +0.38%  '
$ws.Range("E34").Value = '  -4.65%  '
$ws.Range("E35").Value = '  -2.66%  '
$ws.Range("E36").Value = '  -2.23%  '
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.67'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.14%  '
$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.30'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0161'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.95'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.31%  '
$ws.Range("E41").Value = '  -0.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("E43").Value = '  -1.83%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0461'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.13%  '
$ws.Range("E45").Value = '  +2.21%  '
$ws.Range("E46").Value = '  -2.09%  '
$ws.Range("D47").Value = '1.700.74'
$ws.Range("E47").Value = '  -0.56%  '
$ws.Range("E48").Value = '  -5.71%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '85.39'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '43.18'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.31%  '
$ws.Range("E51").Value = '  +0.69%  '
